$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 holds the newly-added record for IP 122.180.21.165 (A3/B3 already populated).
$ws.Range("C3").Value = "Malicious"

$ws.Range("D3").Value = "https://www.virustotal.com/gui/ip-address/122.180.21.165/detection"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://www.virustotal.com/gui/ip-address/122.180.21.165/detection")
$ws.Range("D3").Style = "Hyperlink"

$ws.Range("E3").Value = "{'harmless': 57, 'malicious': 11, 'suspicious': 1, 'undetected': 21, 'timeout': 0}"
$ws.Range("F3").Value = "India"

$ws.Range("G3").Value = 45328.50989583333
$ws.Range("G3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("H3").Value = 45340.6105787037
$ws.Range("H3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("I3").Value = "Bharti Airtel Ltd., Telemedia Services"
